$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.419.70"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.79%  "

$ws.Range("D3").Value = "'2.282.75"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.05%  "

$ws.Range("E4").Value = "  +0.20%  "

$ws.Range("D5").Value = "'322.92"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.80%  "

$ws.Range("D6").Value = "'102.38"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.57%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D9").Value = "'0.607"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.68%  "

$ws.Range("D10").Value = "'39.83"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.16%  "

$ws.Range("D11").Value = "'0.0906"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.10%  "

$ws.Range("D12").Value = "'8.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.96%  "

$ws.Range("D13").Value = "'0.106"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.53%  "

$ws.Range("D14").Value = "'0.966"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.45%  "

$ws.Range("D15").Value = "'15.11"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.63%  "

$ws.Range("D16").Value = "'2.632.65"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.85%  "

$ws.Range("D17").Value = "'2.277.09"
$ws.Range("D17").Style = "Normal"

$ws.Range("D18").Value = "'42.385.97"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.47%  "

$ws.Range("D19").Value = "'7.38"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.04%  "

$ws.Range("E20").Value = "  -0.69%  "

$ws.Range("D21").Value = "'12.79"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +28.09%  "

$ws.Range("E22").Value = "  +1.37%  "

$ws.Range("D23").Value = "'72.99"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "

$ws.Range("D24").Value = "'267.33"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -4.83%  "

$ws.Range("E25").Value = "  -3.11%  "

$ws.Range("E26").Value = "  -0.05%  "

$ws.Range("D27").Value = "'10.85"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.24%  "

$ws.Range("E28").Value = "  +4.17%  "

$ws.Range("D29").Value = "'22.39"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -3.81%  "

$ws.Range("D30").Value = "'37.91"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +6.70%  "

$ws.Range("D31").Value = "'164.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.87%  "

$ws.Range("D32").Value = "'6.10"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.03%  "

$ws.Range("D33").Value = "'0.0878"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.73%  "

$ws.Range("E34").Value = "  +0.50%  "

$ws.Range("E35").Value = "  -13.12%  "

$ws.Range("E36").Value = "  -4.10%  "

$ws.Range("E37").Value = "  -1.69%  "

$ws.Range("E38").Value = "  +0.15%  "

$ws.Range("E39").Value = "  +1.26%  "

$ws.Range("D40").Value = "'2.74"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.87%  "

$ws.Range("D41").Value = "'1.51"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.15%  "

$ws.Range("D42").Value = "'69.39"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.13%  "

$ws.Range("E43").Value = "  +0.00%  "

$ws.Range("E44").Value = "  -1.27%  "

$ws.Range("D45").Value = "'12.39"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.89%  "

$ws.Range("D46").Value = "'90.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -10.80%  "

$ws.Range("D47").Value = "'113.34"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -3.40%  "

$ws.Range("D48").Value = "'79.60"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.98%  "

$ws.Range("E49").Value = "  -2.74%  "

$ws.Range("D50").Value = "'5.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.31%  "

$ws.Range("D51").Value = "'1.594.23"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.06%  "
